$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 26 ("RM 232") - rows below shift up
$ws.Rows.Item(26).Delete()

# After the above delete, former row 28 ("SC 92") is now row 27 - delete it too
$ws.Rows.Item(27).Delete()

# Former row 34 ("SC 193") is now row 32; clear its B value (was -19.9, now blank)
$ws.Cells.Item(32, 2).ClearContents()
